$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (id 12): Dancer
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Dancer"
$ws.Range("C13").Value = "agi"
$ws.Range("D13").Value = "dex"
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = 2
$ws.Range("M13").Value = 0.1

# Row 14 (id 13): Cleric
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Cleric"
$ws.Range("C14").Value = "chr"
$ws.Range("D14").Value = "dex"
$ws.Range("L14").Value = 0.15
$ws.Range("M14").Value = 0.03
$ws.Range("P14").Value = "Fighter"
$ws.Range("Q14").Value = "Prophet"
$ws.Range("R14").Value = 25
$ws.Range("S14").Value = 50

# Row 15 (id 14): Gunslinger
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Gunslinger"
$ws.Range("C15").Value = "agi"
$ws.Range("D15").Value = "dex"
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 2
$ws.Range("L15").Value = 0.12
$ws.Range("M15").Value = 0.08
$ws.Range("O15").Value = 0.04

# Row 16 (id 15): Book Binder
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Book Binder"
$ws.Range("C16").Value = "str"
$ws.Range("D16").Value = "dex"
$ws.Range("M16").Value = 0.15
$ws.Range("P16").Value = "Heretic"
$ws.Range("Q16").Value = "Arcane Alchemist"
$ws.Range("R16").Value = 30
$ws.Range("S16").Value = 50
